# Auto-generated script applying scheduled market-data refresh to Sheets/Marilith_Profits workbook
# Updates cached price/profit figures (columns H-N) on affected leve rows across multiple sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value2 = 144.375
$ws.Range("I33").Value2 = 144.375
$ws.Range("K33").Value2 = 144.375
$ws.Range("M33").Value2 = 84.625
$ws.Range("H62").Value2 = 2533
$ws.Range("I62").Value2 = 2533
$ws.Range("K62").Value2 = 2533
$ws.Range("M62").Value2 = -1909
$ws.Range("H65").Value2 = 2533
$ws.Range("I65").Value2 = 2533
$ws.Range("K65").Value2 = 12665
$ws.Range("M65").Value2 = -9545
$ws.Range("H101").Value2 = 375.42856
$ws.Range("I101").Value2 = 368.6
$ws.Range("J101").Value2 = 392.5
$ws.Range("K101").Value2 = 1105.8
$ws.Range("L101").Value2 = 1177.5
$ws.Range("M101").Value2 = 516.1999999999998
$ws.Range("N101").Value2 = -4421.5
$ws.Range("H137").Value2 = 0
$ws.Range("I137").Value2 = 0
$ws.Range("K137").Value2 = 0
$ws.Range("M137").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value2 = 10000
$ws.Range("I37").Value2 = 10000
$ws.Range("K37").Value2 = 10000
$ws.Range("M37").Value2 = -9727
$ws.Range("H61").Value2 = 2765.4167
$ws.Range("I61").Value2 = 2606.889
$ws.Range("J61").Value2 = 3241
$ws.Range("K61").Value2 = 2606.889
$ws.Range("L61").Value2 = 3241
$ws.Range("M61").Value2 = -2394.889
$ws.Range("N61").Value2 = -3665
$ws.Range("H74").Value2 = 2597.5386
$ws.Range("I74").Value2 = 2597.5386
$ws.Range("J74").Value2 = 0
$ws.Range("K74").Value2 = 2597.5386
$ws.Range("L74").Value2 = 0
$ws.Range("M74").Value2 = -1723.5386
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value2 = 2597.5386
$ws.Range("I77").Value2 = 2597.5386
$ws.Range("J77").Value2 = 0
$ws.Range("K77").Value2 = 12987.693
$ws.Range("L77").Value2 = 0
$ws.Range("M77").Value2 = -8619.692999999999
$ws.Range("N77").ClearContents()
$ws.Range("H122").Value2 = 2656.1875
$ws.Range("I122").Value2 = 2628.5715
$ws.Range("J122").Value2 = 2849.5
$ws.Range("K122").Value2 = 7885.7145
$ws.Range("L122").Value2 = 8548.5
$ws.Range("M122").Value2 = -5435.7145
$ws.Range("N122").Value2 = -13448.5
$ws.Range("H132").Value2 = 1196.1333
$ws.Range("I132").Value2 = 1196.1333
$ws.Range("K132").Value2 = 3588.3999
$ws.Range("M132").Value2 = -1058.3999
$ws.Range("H136").Value2 = 2765.4167
$ws.Range("I136").Value2 = 2606.889
$ws.Range("J136").Value2 = 3241
$ws.Range("K136").Value2 = 7820.667
$ws.Range("L136").Value2 = 9723
$ws.Range("M136").Value2 = -5270.667
$ws.Range("N136").Value2 = -14823

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value2 = 3000
$ws.Range("I94").Value2 = 3000
$ws.Range("J94").Value2 = 0
$ws.Range("K94").Value2 = 3000
$ws.Range("L94").Value2 = 0
$ws.Range("M94").Value2 = -2549
$ws.Range("N94").ClearContents()
$ws.Range("H99").Value2 = 4143.857
$ws.Range("I99").Value2 = 3801.8
$ws.Range("J99").Value2 = 4999
$ws.Range("K99").Value2 = 3801.8
$ws.Range("L99").Value2 = 4999
$ws.Range("M99").Value2 = -2303.8
$ws.Range("N99").Value2 = -7995
$ws.Range("H116").Value2 = 0
$ws.Range("I116").Value2 = 0
$ws.Range("K116").Value2 = 0
$ws.Range("M116").ClearContents()
$ws.Range("H134").Value2 = 4875.7427
$ws.Range("I134").Value2 = 5974.5264
$ws.Range("J134").Value2 = 3570.9375
$ws.Range("K134").Value2 = 17923.5792
$ws.Range("L134").Value2 = 10712.8125
$ws.Range("M134").Value2 = -15388.5792
$ws.Range("N134").Value2 = -15782.8125

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H53").Value2 = 73333.336
$ws.Range("J53").Value2 = 73333.336
$ws.Range("L53").Value2 = 73333.336
$ws.Range("N53").Value2 = -74547.336
$ws.Range("H107").Value2 = 569.8946999999999
$ws.Range("I107").Value2 = 594.53845
$ws.Range("K107").Value2 = 594.53845
$ws.Range("M107").Value2 = 1325.46155

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value2 = 100001
$ws.Range("J37").Value2 = 100001
$ws.Range("L37").Value2 = 300003
$ws.Range("N37").Value2 = -300227
$ws.Range("H74").Value2 = 222853.28
$ws.Range("J74").Value2 = 254993.33
$ws.Range("L74").Value2 = 764979.99
$ws.Range("N74").Value2 = -767101.99
$ws.Range("H75").Value2 = 507.5
$ws.Range("J75").Value2 = 507.5
$ws.Range("L75").Value2 = 1522.5
$ws.Range("N75").Value2 = -3518.5
$ws.Range("H77").Value2 = 222853.28
$ws.Range("J77").Value2 = 254993.33
$ws.Range("L77").Value2 = 2294939.97
$ws.Range("N77").Value2 = -2305547.97
$ws.Range("H78").Value2 = 507.5
$ws.Range("J78").Value2 = 507.5
$ws.Range("L78").Value2 = 4567.5
$ws.Range("N78").Value2 = -14551.5
$ws.Range("H131").Value2 = 1153.75
$ws.Range("I131").Value2 = 940
$ws.Range("J131").Value2 = 1795
$ws.Range("K131").Value2 = 2820
$ws.Range("L131").Value2 = 5385
$ws.Range("M131").Value2 = 2220
$ws.Range("N131").Value2 = -15465
$ws.Range("H132").Value2 = 1898.3334
$ws.Range("J132").Value2 = 1898.3334
$ws.Range("L132").Value2 = 17085.0006
$ws.Range("N132").Value2 = -22145.0006

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value2 = 330
$ws.Range("I55").Value2 = 303
$ws.Range("J55").Value2 = 362.4
$ws.Range("K55").Value2 = 303
$ws.Range("L55").Value2 = 362.4
$ws.Range("M55").Value2 = -130
$ws.Range("N55").Value2 = -708.4
$ws.Range("H63").Value2 = 69085
$ws.Range("J63").Value2 = 69085
$ws.Range("L63").Value2 = 69085
$ws.Range("N63").Value2 = -70583
$ws.Range("H66").Value2 = 69085
$ws.Range("J66").Value2 = 69085
$ws.Range("L66").Value2 = 207255
$ws.Range("N66").Value2 = -214743
$ws.Range("H100").Value2 = 4915.6665
$ws.Range("I100").Value2 = 4498.5
$ws.Range("K100").Value2 = 4498.5
$ws.Range("M100").Value2 = -3957.5
$ws.Range("H139").Value2 = 49999
$ws.Range("J139").Value2 = 49999
$ws.Range("L139").Value2 = 49999
$ws.Range("N139").Value2 = -60279

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value2 = 30381
$ws.Range("I16").Value2 = 30381
$ws.Range("K16").Value2 = 30381
$ws.Range("M16").Value2 = -30089
$ws.Range("H54").Value2 = 1013333.3
$ws.Range("I54").Value2 = 0
$ws.Range("J54").Value2 = 1013333.3
$ws.Range("K54").Value2 = 0
$ws.Range("L54").Value2 = 1013333.3
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value2 = -1014373.3
$ws.Range("H107").Value2 = 1143.2
$ws.Range("I107").Value2 = 804.36365
$ws.Range("J107").Value2 = 2075
$ws.Range("K107").Value2 = 2413.09095
$ws.Range("L107").Value2 = 6225
$ws.Range("M107").Value2 = -493.0909499999998
$ws.Range("N107").Value2 = -10065
$ws.Range("H132").Value2 = 2897.5
$ws.Range("I132").Value2 = 799.3333
$ws.Range("K132").Value2 = 2397.9999
$ws.Range("M132").Value2 = 132.0001000000002
